$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 117.18
$ws.Range("I15").Value = 117.18
$ws.Range("K15").Value = 351.54
$ws.Range("M15").Value = -182.54
# Row 21
$ws.Range("H21").Value = 8142.857
$ws.Range("J21").Value = 8142.857
$ws.Range("L21").Value = 8142.857
$ws.Range("N21").Value = -9078.857
# Row 23
$ws.Range("H23").Value = 8142.857
$ws.Range("J23").Value = 8142.857
$ws.Range("L23").Value = 8142.857
$ws.Range("N23").Value = -8610.857
# Row 62
$ws.Range("H62").Value = 4666.6665
$ws.Range("I62").Value = 4250
$ws.Range("K62").Value = 4250
$ws.Range("M62").Value = -3626
# Row 65
$ws.Range("H65").Value = 4666.6665
$ws.Range("I65").Value = 4250
$ws.Range("K65").Value = 21250
$ws.Range("M65").Value = -18130
# Row 76
$ws.Range("H76").Value = 82432.42999999999
$ws.Range("I76").Value = 108625.336
$ws.Range("J76").Value = 3853.7144
$ws.Range("K76").Value = 108625.336
$ws.Range("L76").Value = 3853.7144
$ws.Range("M76").Value = -108310.336
$ws.Range("N76").Value = -4483.7144
# Row 79
$ws.Range("H79").Value = 82432.42999999999
$ws.Range("I79").Value = 108625.336
$ws.Range("J79").Value = 3853.7144
$ws.Range("K79").Value = 108625.336
$ws.Range("L79").Value = 3853.7144
$ws.Range("M79").Value = -107533.336
$ws.Range("N79").Value = -6037.7144
# Row 92
$ws.Range("H92").Value = 589.63635
$ws.Range("I92").Value = 558.8570999999999
$ws.Range("J92").Value = 643.5
$ws.Range("K92").Value = 558.8570999999999
$ws.Range("L92").Value = 643.5
$ws.Range("M92").Value = 689.1429000000001
$ws.Range("N92").Value = -3139.5
# Row 132
$ws.Range("H132").Value = 3345.0725
$ws.Range("I132").Value = 3300.111
$ws.Range("J132").Value = 3506.9333
$ws.Range("K132").Value = 9900.332999999999
$ws.Range("L132").Value = 10520.7999
$ws.Range("M132").Value = -7370.332999999999
$ws.Range("N132").Value = -15580.7999
# Row 137
$ws.Range("H137").Value = 26103.85
$ws.Range("I137").Value = 844.25
$ws.Range("J137").Value = 85042.914
$ws.Range("K137").Value = 2532.75
$ws.Range("L137").Value = 255128.742
$ws.Range("M137").Value = 17.25
$ws.Range("N137").Value = -260228.742

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1359.8966
$ws.Range("I45").Value = 1427.409
$ws.Range("K45").Value = 1427.409
$ws.Range("M45").Value = -1050.409
# Row 61
$ws.Range("H61").Value = 1999.3226
$ws.Range("I61").Value = 1894.2609
$ws.Range("J61").Value = 2301.375
$ws.Range("K61").Value = 1894.2609
$ws.Range("L61").Value = 2301.375
$ws.Range("M61").Value = -1682.2609
$ws.Range("N61").Value = -2725.375
# Row 74
$ws.Range("H74").Value = 229511.1
$ws.Range("I74").Value = 2269.4119
$ws.Range("K74").Value = 2269.4119
$ws.Range("M74").Value = -1395.4119
# Row 77
$ws.Range("H77").Value = 229511.1
$ws.Range("I77").Value = 2269.4119
$ws.Range("K77").Value = 11347.0595
$ws.Range("M77").Value = -6979.059499999999
# Row 132
$ws.Range("H132").Value = 27799.46
$ws.Range("I132").Value = 32849.5
$ws.Range("J132").Value = 4713.5713
$ws.Range("K132").Value = 98548.5
$ws.Range("L132").Value = 14140.7139
$ws.Range("M132").Value = -96018.5
$ws.Range("N132").Value = -19200.7139
# Row 136
$ws.Range("H136").Value = 1999.3226
$ws.Range("I136").Value = 1894.2609
$ws.Range("J136").Value = 2301.375
$ws.Range("K136").Value = 5682.7827
$ws.Range("L136").Value = 6904.125
$ws.Range("M136").Value = -3132.7827
$ws.Range("N136").Value = -12004.125

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6898581.5
$ws.Range("I86").Value = 8335244.5
$ws.Range("J86").Value = 2597.8
$ws.Range("K86").Value = 8335244.5
$ws.Range("L86").Value = 2597.8
$ws.Range("M86").Value = -8334121.5
$ws.Range("N86").Value = -4843.8
# Row 89
$ws.Range("H89").Value = 6898581.5
$ws.Range("I89").Value = 8335244.5
$ws.Range("J89").Value = 2597.8
$ws.Range("K89").Value = 41676222.5
$ws.Range("L89").Value = 12989
$ws.Range("M89").Value = -41670606.5
$ws.Range("N89").Value = -24221
# Row 107
$ws.Range("H107").Value = 4233.147
$ws.Range("I107").Value = 5200.48
$ws.Range("J107").Value = 1546.1111
$ws.Range("K107").Value = 5200.48
$ws.Range("L107").Value = 1546.1111
$ws.Range("M107").Value = -3280.48
$ws.Range("N107").Value = -5386.1111
# Row 134
$ws.Range("H134").Value = 590446.9
$ws.Range("I134").Value = 1111994.1
$ws.Range("J134").Value = 3706.25
$ws.Range("K134").Value = 3335982.3
$ws.Range("L134").Value = 11118.75
$ws.Range("M134").Value = -3333447.3
$ws.Range("N134").Value = -16188.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1939.7797
$ws.Range("I31").Value = 1177.5
$ws.Range("J31").Value = 2540.3635
$ws.Range("K31").Value = 1177.5
$ws.Range("L31").Value = 2540.3635
$ws.Range("M31").Value = -882.5
$ws.Range("N31").Value = -3130.3635
# Row 34
$ws.Range("H34").Value = 1939.7797
$ws.Range("I34").Value = 1177.5
$ws.Range("J34").Value = 2540.3635
$ws.Range("K34").Value = 1177.5
$ws.Range("L34").Value = 2540.3635
$ws.Range("M34").Value = -975.5
$ws.Range("N34").Value = -2944.3635
# Row 58
$ws.Range("H58").Value = 2857.3845
$ws.Range("I58").Value = 1113.8667
$ws.Range("J58").Value = 3947.0833
$ws.Range("K58").Value = 1113.8667
$ws.Range("L58").Value = 3947.0833
$ws.Range("M58").Value = -910.8667
$ws.Range("N58").Value = -4353.0833
# Row 86
$ws.Range("H86").Value = 252113
$ws.Range("I86").Value = 335284
$ws.Range("J86").Value = 2600
$ws.Range("K86").Value = 335284
$ws.Range("L86").Value = 2600
$ws.Range("M86").Value = -334161
$ws.Range("N86").Value = -4846
# Row 89
$ws.Range("H89").Value = 252113
$ws.Range("I89").Value = 335284
$ws.Range("J89").Value = 2600
$ws.Range("K89").Value = 1676420
$ws.Range("L89").Value = 13000
$ws.Range("M89").Value = -1670804
$ws.Range("N89").Value = -24232
# Row 107
$ws.Range("H107").Value = 2641.4285
$ws.Range("I107").Value = 2298.2
$ws.Range("J107").Value = 3499.5
$ws.Range("K107").Value = 2298.2
$ws.Range("L107").Value = 3499.5
$ws.Range("M107").Value = -378.1999999999998
$ws.Range("N107").Value = -7339.5
# Row 132
$ws.Range("H132").Value = 2408.25
$ws.Range("I132").Value = 2267.5625
$ws.Range("J132").Value = 2689.625
$ws.Range("K132").Value = 6802.6875
$ws.Range("L132").Value = 8068.875
$ws.Range("M132").Value = -4272.6875
$ws.Range("N132").Value = -13128.875
# Row 134
$ws.Range("H134").Value = 3527.1628
$ws.Range("I134").Value = 3816.516
$ws.Range("J134").Value = 2779.6667
$ws.Range("K134").Value = 11449.548
$ws.Range("L134").Value = 8339.000100000001
$ws.Range("M134").Value = -8914.548000000001
$ws.Range("N134").Value = -13409.0001
# Row 136
$ws.Range("H136").Value = 2857.3845
$ws.Range("I136").Value = 1113.8667
$ws.Range("J136").Value = 3947.0833
$ws.Range("K136").Value = 3341.6001
$ws.Range("L136").Value = 11841.2499
$ws.Range("M136").Value = -791.6001000000001
$ws.Range("N136").Value = -16941.2499

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 3945.2222
$ws.Range("I33").Value = 170.85715
$ws.Range("J33").Value = 5266.25
$ws.Range("K33").Value = 1025.1429
$ws.Range("L33").Value = 31597.5
$ws.Range("M33").Value = -742.1428999999998
$ws.Range("N33").Value = -32163.5
# Row 68
$ws.Range("H68").Value = 1411.25
$ws.Range("J68").Value = 1548.3334
$ws.Range("L68").Value = 4645.0002
$ws.Range("N68").Value = -6267.0002
# Row 71
$ws.Range("H71").Value = 1411.25
$ws.Range("J71").Value = 1548.3334
$ws.Range("L71").Value = 13935.0006
$ws.Range("N71").Value = -22047.0006
# Row 113
$ws.Range("H113").Value = 1137.9246
$ws.Range("I113").Value = 441.8
$ws.Range("J113").Value = 1412.7106
$ws.Range("K113").Value = 1325.4
$ws.Range("L113").Value = 4238.1318
$ws.Range("M113").Value = 844.5999999999999
$ws.Range("N113").Value = -8578.131799999999
# Row 131
$ws.Range("H131").Value = 820.2414
$ws.Range("I131").Value = 408.1
$ws.Range("J131").Value = 1037.1578
$ws.Range("K131").Value = 1224.3
$ws.Range("L131").Value = 3111.4734
$ws.Range("M131").Value = 3815.7
$ws.Range("N131").Value = -13191.4734

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3129.0454
$ws.Range("I126").Value = 1938.3846
$ws.Range("K126").Value = 5815.1538
$ws.Range("M126").Value = -3345.1538
# Row 132
$ws.Range("H132").Value = 2378.9395
$ws.Range("I132").Value = 1710.7778
$ws.Range("J132").Value = 3180.7334
$ws.Range("K132").Value = 5132.3334
$ws.Range("L132").Value = 9542.200199999999
$ws.Range("M132").Value = -2602.3334
$ws.Range("N132").Value = -14602.2002

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 329.4
$ws.Range("I16").Value = 220.84616
$ws.Range("J16").Value = 531
$ws.Range("K16").Value = 220.84616
$ws.Range("L16").Value = 531
$ws.Range("M16").Value = -50.84616
$ws.Range("N16").Value = -871
# Row 132
$ws.Range("H132").Value = 16138298
$ws.Range("I132").Value = 31264606
$ws.Range("J132").Value = 3569.2
$ws.Range("K132").Value = 93793818
$ws.Range("L132").Value = 10707.6
$ws.Range("M132").Value = -93791288
$ws.Range("N132").Value = -15767.6

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1656.25
$ws.Range("I132").Value = 897.5517
$ws.Range("J132").Value = 3123.0667
$ws.Range("K132").Value = 2692.6551
$ws.Range("L132").Value = 9369.2001
$ws.Range("M132").Value = -162.6550999999999
$ws.Range("N132").Value = -14429.2001
# Row 136
$ws.Range("H136").Value = 8466.888999999999
$ws.Range("I136").Value = 3013
$ws.Range("J136").Value = 12830
$ws.Range("K136").Value = 9039
$ws.Range("L136").Value = 38490
$ws.Range("M136").Value = -6489
$ws.Range("N136").Value = -43590
